$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update D23: 2.5 -> 3
$ws.Range("D23").Value = 3

# Add D24: new value 3
$ws.Range("D24").Value = 3

# Update C25: 2 -> 4
$ws.Range("C25").Value = 4

# Add D25: new value 3
$ws.Range("D25").Value = 3

# Update the active selection to C39
$ws.Range("C39").Select()
